# Add four new "charge" columns (IMPORT_CHARGES, EXPORT_CHARGES, PRE_CARRIAGE,
# ON_CARRIAGE) to the hubs worksheet, right after the existing PHOTO column,
# matching the "pre merge with obi" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing header row is A1:I1 (STATUS .. PHOTO). New headers go in J1:M1.
$ws.Range("J1").Value = "IMPORT_CHARGES"
$ws.Range("K1").Value = "EXPORT_CHARGES"
$ws.Range("L1").Value = "PRE_CARRIAGE"
$ws.Range("M1").Value = "ON_CARRIAGE"

# Give the new header cells the same (bold header) formatting as the rest of
# row 1 by copying the format from the last existing header cell (I1).
$ws.Range("I1").Copy() | Out-Null
$ws.Range("J1:M1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
